$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.800.04'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '3.803.65'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').Value = "'596.70"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = "'167.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').Value = '3.801.11'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = "'0.161"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('D12').Value = "'0.450"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('E13').Value = '  -2.34%  '
$ws.Range('D14').Value = "'36.09"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').Value = '4.440.97'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = '3.777.11'
$ws.Range('E16').Value = '  -2.32%  '
$ws.Range('D17').Value = "'18.59"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.60%  '
$ws.Range('D18').Value = '67.781.19'
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').Value = "'7.11"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = "'461.70"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').Value = "'9.89"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.77%  '
$ws.Range('D23').Value = "'0.702"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').Value = "'0.0000154"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  -0.82%  '
$ws.Range('D26').Value = "'12.06"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('E27').Value = '  -2.74%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = "'10.01"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('D30').Value = '3.948.85'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('D34').Value = "'29.70"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = "'9.07"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = "'0.0999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = "'3.38"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = "'48.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').Value = "'44.08"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('D47').Value = "'150.71"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.92%  '
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').Value = "'392.37"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('D50').Value = "'1.83"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('D51').Value = "'26.36"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.95%  '
